# Query Fixed for TC05,9,12,26,28,30,36,39,40,42,46 and moved to Breed
#
# The "StatQuery" column (column C, rows 2-4 on the "startup" sheet) held an
# outdated Cypher query used to compute summary statistics. Replace it with
# the corrected/expanded query (now also returning aliquot counts and using
# the parameterized filter style), and move the active selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("startup")

$newStatQuery = " MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n" +
    "      WHERE (size([]) = 0 OR s.clinical_study_designation IN [])`n" +
    "        AND (s.study_disposition = 'Unrestricted')`n" +
    "        AND (size([]) = 0 OR s.clinical_study_type IN [])`n" +
    "        AND (size(['Saint Bernard']) = 0 OR demo.breed IN ['Saint Bernard'])`n" +
    "        AND (size([]) = 0 OR demo.sex IN [])`n" +
    "        AND (size([]) = 0 OR demo.neutered_indicator IN [])`n" +
    "        AND (size([]) = 0 OR diag.disease_term IN [])`n" +
    "        AND (size([]) = 0 OR diag.primary_disease_site IN [])`n" +
    "        AND (size([]) = 0 OR diag.stage_of_disease IN [])`n" +
    "        AND (size([]) = 0 OR diag.best_response IN [])`n" +
    "    OPTIONAL MATCH (c)-->(co:cohort)`n" +
    "    OPTIONAL MATCH (f:file)-[*]->(c)`n" +
    "    OPTIONAL MATCH (f)-->(parent)`n" +
    "    OPTIONAL MATCH (samp:sample)-->(c)`n" +
    "    OPTIONAL MATCH (samp)<--(al:aliquot)`n" +
    "    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al`n" +
    "      WHERE (size([]) = 0 OR samp.summarized_sample_type IN [])`n" +
    "        AND (size([]) = 0 OR samp.specific_sample_pathology IN [])`n" +
    "        AND (size([]) = 0 OR samp.sample_site IN [])`n" +
    "        AND (size([]) = 0 OR head(labels(parent)) IN [])`n" +
    "        AND (size([]) = 0 OR f.file_type IN [])`n" +
    "        AND (size([]) = 0 OR f.file_format IN [])`n" +
    "    WITH c.case_id AS case_id,`n" +
    "         s.clinical_study_designation AS study_code,`n" +
    "         s.clinical_study_type AS study_type,`n" +
    "         co.cohort_description AS cohort,`n" +
    "         demo.breed AS breed,`n" +
    "         diag.disease_term AS diagnosis,`n" +
    "         diag.stage_of_disease AS stage_of_disease,`n" +
    "         diag.primary_disease_site AS disease_site,`n" +
    "         demo.patient_age_at_enrollment AS age,`n" +
    "         demo.sex AS sex,`n" +
    "         demo.neutered_indicator AS neutered_status,`n" +
    "         demo.weight AS weight,`n" +
    "         diag.best_response AS response_to_treatment,`n" +
    "         samp.sample_id AS sample_id,`n" +
    "         f.uuid AS file_id,`n" +
    "         al`n" +
    "    RETURN`n" +
    "COUNT(DISTINCT file_id) as number_of_files,`n" +
    "COUNT(DISTINCT sample_id) as number_of_sample,`n" +
    "COUNT(DISTINCT case_id) as number_of_cases,`n" +
    "COUNT(DISTINCT study_code) as number_of_study,`n" +
    "COUNT(DISTINCT al) as number_of_aliquot"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

$ws.Rows.Item(2).RowHeight = 409.6
$ws.Rows.Item(3).RowHeight = 409.6
$ws.Rows.Item(4).RowHeight = 409.6

$ws.Range("B4").Select()
